# "Add zip-excel export writer" - the course-export template's placeholder
# for the applicant's origin column is switched from the short-name token
# to the full-name token.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kursliste")

$ws.Range("D2").Value = "applicant.origin.name"
$ws.Range("D2").Select()
